$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule (cell E8) to reflect the git update
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell selection left by the edit
$ws.Range("E8").Select()
